# Update product list on "Tabelle1" (sheet1):
#  - fix "Chedar" -> "Cheddar" spelling in the Hamburger ingredients text
#  - Beilagen (sides): drop "Potatoes", keep "Gemischter Salat"
#  - Extension (toppings): add a new "Tomate" row, shift rows down one,
#    rename the category label on the new last row to "Extension " (trailing space)
#  - Saucen (sauces): rename "Majonnaise" -> "Mayonnaise" and add new
#    "Cocktail" / "Curry" sauces, shifting rows down one
#  - view state: scroll down a bit and select B24

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Hamburger ingredients spelling fix ---
$ws.Range("D5").Value = "Weichbrötchen, Rindfleischburger, Cheddar, Gurke"

# --- Beilagen (sides): remove "Potatoes" row, "Gemischter Salat" moves up ---
$ws.Range("B13").Value = "Gemischter Salat"
$ws.Range("A14").ClearContents()
$ws.Range("B14").ClearContents()

# --- Extension (toppings): insert new row 16, shift rows 17-21 down ---
$ws.Range("A16").Value = "Extension"
$ws.Range("B16").Value = "Speck Streifen"
$ws.Range("B17").Value = "Käse"
$ws.Range("B18").Value = "Rindfleischburger extra"
$ws.Range("B19").Value = "Gurke"
$ws.Range("B20").Value = "Chili"
$ws.Range("A21").Value = "Extension "
$ws.Range("B21").Value = "Tomate"

# --- Saucen (sauces): insert new row 24, shift rows 25-29 down, add row 30 ---
$ws.Range("A24").Value = "Saucen"
$ws.Range("B24").Value = "Senf"
$ws.Range("B25").Value = "Ketchup"
$ws.Range("B26").Value = "Barbecue"
$ws.Range("B27").Value = "Mayonnaise"
$ws.Range("B28").Value = "Hot & Spicy"
$ws.Range("A29").Value = "Saucen"
$ws.Range("B29").Value = "Cocktail"
$ws.Range("A30").Value = "Saucen"
$ws.Range("B30").Value = "Curry"

# --- view state: scroll the window and select B24 ---
$excel.ActiveWindow.ScrollRow = 3
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("B24").Select()
